$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "zmin"
$ws.Range("K1").Value = "zmax"

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 10).Value = 0.2
    $ws.Cells.Item($r, 11).Value = 1
}

$ws.Range("K11").Select()
